# Merged Lavanya's Code on 14th september
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Description field): value text got a "1234" suffix.
$ws.Range("B5").Value = "Description1234"

# Rows 8-9 used to describe the organisation dropdown/tree locators;
# they were replaced with the username/search-button locators.
$ws.Range("A8").Value = "text"
$ws.Range("B8").Value = "Username"
$ws.Range("C8").Value = "id"
$ws.Range("D8").Value = "__ns2036994902_name"

$ws.Range("A9").Value = "button"
$ws.Range("B9").Value = "Search"
$ws.Range("C9").Value = "id"
$ws.Range("D9").Value = "__ns2036994902_searchFiltersSearchBtnText"

# New locator rows appended for the provisioning-profile / GHLR workflow.
$ws.Range("A10").Value = "inputlist"
$ws.Range("B10").Value = "Organisation"
$ws.Range("C10").Value = "class"
$ws.Range("D10").Value = "fieldBlock__field fieldBlock__field--withClean"

$ws.Range("A11").Value = "table"
$ws.Range("B11").Value = "ServiceProf"
$ws.Range("C11").Value = "xpath"
$ws.Range("D11").Value = "//*[@class='globaltable dataTable']"

$ws.Range("A12").Value = "text"
$ws.Range("B12").Value = "ServingNetworkCode"
$ws.Range("C12").Value = "id"
$ws.Range("D12").Value = "__ns2036994902_searchCode"

$ws.Range("A13").Value = "text"
$ws.Range("B13").Value = "GHLRtemp"
$ws.Range("C13").Value = "id"
$ws.Range("D13").Value = "__ns2036994902_searchGhlrNameList"

$ws.Range("A14").Value = "text"
$ws.Range("B14").Value = "Service Profile"
$ws.Range("C14").Value = "id"
$ws.Range("D14").Value = "__ns2036994902_serviceProfileInput"

$ws.Range("A15").Value = "button"
$ws.Range("B15").Value = "ConfigEdit"
$ws.Range("C15").Value = "id"
$ws.Range("D15").Value = "__ns2036994902_configurationAreaEditBtn"

$ws.Range("A16").Value = "text"
$ws.Range("B16").Value = "Provisioning profile"
$ws.Range("C16").Value = "id"
$ws.Range("D16").Value = "__ns2036994902_searchProvisioningProfile"

$ws.Range("A17").Value = "button"
$ws.Range("B17").Value = "Create provisioning profile"
$ws.Range("C17").Value = "id"
$ws.Range("D17").Value = "__ns2036994902_createBtn"

$ws.Range("A18").Value = "date"
$ws.Range("B18").Value = "Withdraw Date (UTC)"
$ws.Range("C18").Value = "id"
$ws.Range("D18").Value = "__ns2036993941_withdrawDateCalendar"

# Row 10 is highlighted yellow to flag it for review.
$ws.Range("A10:D10").Interior.Color = 65535

# Leave the selection on B5, matching where the edit was made.
$ws.Range("B5").Select()
